# daily auto push: 2026-01-29 14:04 UTC
#
# Two new observation rows for 2026/01/29 (Thu) are inserted at the top of
# the "2026/12/29" block (row 719), pushing all subsequent rows down by 2
# and extending the used range from A1:D760 to A1:D762.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 719 onward down by two, preserving all of their existing data.
$ws.Rows("719:720").Insert()

# Row 719: new record for 2026/01/29 (Thu)
$ws.Cells.Item(719, 1).NumberFormat = "@"
$ws.Cells.Item(719, 1).Value = "2026/01/29"
$ws.Cells.Item(719, 2).Value = "木"
$ws.Cells.Item(719, 3).Value = 18
$ws.Cells.Item(719, 4).Value = 19

# Row 720: new record for 2026/01/29 (Thu)
$ws.Cells.Item(720, 1).NumberFormat = "@"
$ws.Cells.Item(720, 1).Value = "2026/01/29"
$ws.Cells.Item(720, 2).Value = "木"
$ws.Cells.Item(720, 3).Value = 22
$ws.Cells.Item(720, 4).Value = 21
